$wb = $excel.ActiveWorkbook

# xlPasteValues constant (used to paste a text "True" value without Excel
# re-interpreting the literal string "True" as a Boolean cell type)
$xlPasteValues = -4163
# xlCenter constant
$xlCenter = -4108

# A helper "True" text cell already present on the StorageUnit sheet
# (stored as a shared string, not a boolean) - used as a copy source so new
# "True" cells keep the same t="s" shared-string representation.
$wsStorage = $wb.Worksheets.Item("StorageUnit")

# --- Line sheet: add new row 2 "line_0-1" ---
$wsLine = $wb.Worksheets.Item("Line")
$wsLine.Range("A2").Value = "line_0-1"
$wsLine.Range("B2").Value = "AC"
$wsLine.Range("C2").Value = "bus 0"
$wsLine.Range("D2").Value = "bus 1"
$wsLine.Range("F2").Value = 0.1
$wsLine.Range("G2").Value = 0.1
$wsLine.Range("A2:G2").HorizontalAlignment = $xlCenter

$wsStorage.Range("D2").Copy()
$wsLine.Range("E2").PasteSpecial($xlPasteValues)
$excel.CutCopyMode = 0
$wsLine.Range("E2").HorizontalAlignment = $xlCenter

$wsLine.Range("G3").Select()

# --- Generator sheet: add new row 3 "diesel" ---
$wsGen = $wb.Worksheets.Item("Generator")
$wsGen.Range("A3").Value = "diesel"
$wsGen.Range("B3").Value = "diesel"
$wsGen.Range("C3").Value = "bus 0"
$wsGen.Range("E3").Value = 0
$wsGen.Range("F3").Value = 1
$wsGen.Range("G3").Value = 100
$wsGen.Range("A3:G3").HorizontalAlignment = $xlCenter

$wsStorage.Range("D2").Copy()
$wsGen.Range("D3").PasteSpecial($xlPasteValues)
$excel.CutCopyMode = 0
$wsGen.Range("D3").HorizontalAlignment = $xlCenter

$wsGen.Range("G4").Select()

# --- Link sheet: remove the battery_link row (row 2) ---
$wsLink = $wb.Worksheets.Item("Link")
$wsLink.Rows.Item(2).Delete()
$wsLink.Rows.Item(2).EntireRow.Select()

# restore the originally active sheet/selection (StorageUnit) so the
# workbook-level active tab is left untouched by this edit
$wsStorage.Activate()
$wsStorage.Range("H3").Select()
